$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before the current row 312, shifting existing
# rows 312:336 down to 315:339 (preserving all their data/formatting).
$ws.Rows("312:314").Insert()

# Static field values shared by every data row in this sheet.
$mercadoId   = 6
$mercado     = "Mercado Mayorista Lo Valledor de Santiago"
$region      = "Metropolitana"
$productoId  = 13
$productoCod = 100112043
$producto    = "Pepino dulce"
$variedad    = "Cultivar IV Región"
$unidad      = "`$/bandeja 18 kilos"
$provincia   = "Provincia de Limarí"
$kilos       = 18
$tipo        = "Hortaliza"

# New weekly data (2023-04-05, serial 45021) for the three inserted rows.
$newRows = @(
    @{ Row = 312; Calidad = "Especial"; Muestras = 580; Minimo = 13000; Maximo = 14000; Medio = 13603; PrecioKg = 756 },
    @{ Row = 313; Calidad = "Primera";  Muestras = 430; Minimo = 11000; Maximo = 12000; Medio = 11605; PrecioKg = 645 },
    @{ Row = 314; Calidad = "Segunda";  Muestras = 320; Minimo = 9000;  Maximo = 10000; Medio = 9531;  PrecioKg = 530 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = 45021
    $ws.Cells.Item($row, 5).Value = $productoId
    $ws.Cells.Item($row, 6).Value = $productoCod
    $ws.Cells.Item($row, 7).Value = $producto
    $ws.Cells.Item($row, 8).Value = $variedad
    $ws.Cells.Item($row, 9).Value = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Muestras
    $ws.Cells.Item($row, 11).Value = $r.Minimo
    $ws.Cells.Item($row, 12).Value = $r.Maximo
    $ws.Cells.Item($row, 13).Value = $r.Medio
    $ws.Cells.Item($row, 14).Value = $unidad
    $ws.Cells.Item($row, 15).Value = $provincia
    $ws.Cells.Item($row, 16).Value = $r.PrecioKg
    $ws.Cells.Item($row, 17).Value = $kilos
    $ws.Cells.Item($row, 18).Value = $tipo
}
